$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("I5").Value = 11
$ws.Range("K5").Value = 2.3
$ws.Range("L5").Value = 9.5
$ws.Range("M5").Value = 1.05
$ws.Range("N5").Value = 11
$ws.Range("O5").Value = 1.3
$ws.Range("P5").Value = 3.4
$ws.Range("Q5").Value = 2.03
$ws.Range("R5").Value = 1.83
$ws.Range("S5").Value = 1.4
$ws.Range("T5").Value = 2.75
$ws.Range("W5").Value = 5.5
$ws.Range("X5").Value = 5.5
$ws.Range("Y5").Value = 9.5
$ws.Range("AA5").Value = 13
$ws.Range("AC5").Value = 8.5
$ws.Range("AF5").Value = 101
$ws.Range("AI5").Value = 41
$ws.Range("AJ5").Value = 29
$ws.Range("AL5").Value = 81
$ws.Range("AM5").Value = 81
$ws.Range("AS5").Value = 201
$ws.Range("AT5").Value = 2.75
$ws.Range("AU5").Value = 11
$ws.Range("AV5").Value = 81
$ws.Range("AZ5").Value = 301

# Row 11
$ws.Range("G11").Value = 1.53
$ws.Range("H11").Value = 3.8
$ws.Range("I11").Value = 7
$ws.Range("M11").Value = 1.08
$ws.Range("N11").Value = 8
$ws.Range("W11").Value = 6
$ws.Range("X11").Value = 6.5
$ws.Range("Z11").Value = 10
$ws.Range("AC11").Value = 8
$ws.Range("AD11").Value = 7.5
$ws.Range("AF11").Value = 67
$ws.Range("AO11").Value = 8
$ws.Range("AP11").Value = 21
$ws.Range("AU11").Value = 9.5
$ws.Range("AV11").Value = 67
$ws.Range("AW11").Value = 7.5

# Row 13
$ws.Range("G13").Value = 1.75
$ws.Range("H13").Value = 3.25
$ws.Range("I13").Value = 4.8
$ws.Range("J13").Value = 2.35
$ws.Range("K13").Value = 2.02
$ws.Range("L13").Value = 4.9
$ws.Range("N13").Value = 8.9
$ws.Range("O13").Value = 1.29
$ws.Range("P13").Value = 3
$ws.Range("Q13").Value = 1.91
$ws.Range("R13").Value = 1.8
$ws.Range("U13").Value = 1.75
$ws.Range("V13").Value = 1.87
$ws.Range("W13").Value = 6.6
$ws.Range("X13").Value = 8
$ws.Range("Y13").Value = 7.9
$ws.Range("Z13").Value = 14.5
$ws.Range("AA13").Value = 14
$ws.Range("AB13").Value = 26
$ws.Range("AC13").Value = 9
$ws.Range("AD13").Value = 6.4
$ws.Range("AE13").Value = 14
$ws.Range("AF13").Value = 65
$ws.Range("AH13").Value = 13.5
$ws.Range("AJ13").Value = 15
$ws.Range("AK13").Value = 90
$ws.Range("AL13").Value = 45
$ws.Range("AM13").Value = 45
$ws.Range("AN13").Value = 3.55
$ws.Range("AO13").Value = 9
$ws.Range("AP13").Value = 18
$ws.Range("AQ13").Value = 32
$ws.Range("AR13").Value = 65
$ws.Range("AS13").Value = 250
$ws.Range("AU13").Value = 7
$ws.Range("AW13").Value = 6.4
$ws.Range("AX13").Value = 27
$ws.Range("AY13").Value = 30
$ws.Range("AZ13").Value = 175
$ws.Range("BA13").Value = 175
$ws.Range("BB13").Value = 400
